$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before column N so that the old N:Q (מספר חוזה, מספר מונה, סוג התעריף, שם בן הזוג)
# shift to R:U, making room for the new card columns.
$ws.Range("N1:Q1").EntireColumn.Insert()

# New header values for the inserted columns N1:Q1
$ws.Range("N1").Value = "שם בעל הכרטיס"
$ws.Range("O1").Value = "מספר זהות/דרכון"
$ws.Range("P1").Value = "מספר כרטיס"
$ws.Range("Q1").Value = "תוקף"

# Clear the new card columns for existing data rows (2-6 for now, rows 4-6 will be deleted below)
$ws.Range("N2:Q6").ClearContents()

# Delete data rows 4, 5 and 6 - only rows 2 and 3 of data remain
$ws.Range("A6:U6").EntireRow.Delete()
$ws.Range("A5:U5").EntireRow.Delete()
$ws.Range("A4:U4").EntireRow.Delete()

# Update row 2 values
$ws.Range("A2").Value = "חני"
$ws.Range("C2").Value = ""

# Update row 3 values
$ws.Range("A3").Value = "מאיר"
$ws.Range("B3").Value = "חנה"
$ws.Range("C3").Value = "כהן"

$wb.Save()
